$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Shift rows 11-14 down to 12-15 (bottom-up) to make room for a new test
#    case row at row 11. Copy formats first, then values, so per-cell styles
#    (which vary per row, e.g. quotePrefix on column B) survive intact.
# ---------------------------------------------------------------------------
$ws.Range("A14:K14").Copy() | Out-Null
$ws.Range("A15:K15").PasteSpecial(-4122) | Out-Null
$ws.Range("A14:K14").Copy() | Out-Null
$ws.Range("A15:K15").PasteSpecial(-4163) | Out-Null

$ws.Range("A13:K13").Copy() | Out-Null
$ws.Range("A14:K14").PasteSpecial(-4122) | Out-Null
$ws.Range("A13:K13").Copy() | Out-Null
$ws.Range("A14:K14").PasteSpecial(-4163) | Out-Null

$ws.Range("A12:K12").Copy() | Out-Null
$ws.Range("A13:K13").PasteSpecial(-4122) | Out-Null
$ws.Range("A12:K12").Copy() | Out-Null
$ws.Range("A13:K13").PasteSpecial(-4163) | Out-Null

$ws.Range("A11:K11").Copy() | Out-Null
$ws.Range("A12:K12").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:K11").Copy() | Out-Null
$ws.Range("A12:K12").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Populate the new row 11 (keeps the formatting that row 11 already had,
#    which is the same style class used by row 10/row 12).
# ---------------------------------------------------------------------------
$ws.Range("A11").Value2 = "TC_Negatice"
$ws.Range("B11").Value2 = "Yes"
$ws.Range("C11").Value2 = "ABCD"
$ws.Range("D11").Value2 = "SP"
$ws.Range("E11").Value2 = "IQ"
$ws.Range("F11").Value2 = "ar"
$ws.Range("G11").Value2 = "mobile"
$ws.Range("H11").Value2 = "android"
$ws.Range("I11").Value2 = 18
$ws.Range("J11").Value2 = 1
$ws.Range("K11").Value2 = "0-100"

# ---------------------------------------------------------------------------
# 3) Global flag change: set column B (RunMode) to "Yes" for every data row.
#    Writing a new value to a quote-prefixed cell (style s=3, used by col B)
#    makes the runtime drop the quotePrefix flag, so re-apply the original
#    per-cell format (copied from B2, which already is "Yes"/style s=3)
#    immediately afterwards without touching the value.
# ---------------------------------------------------------------------------
$bRows = 3..15
foreach ($r in $bRows) {
    $ws.Range("B$r").Value2 = "Yes"
}
$ws.Range("B2").Copy() | Out-Null
foreach ($r in $bRows) {
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Row 2 / column C changes from "home" to "majidtv".
# ---------------------------------------------------------------------------
$ws.Range("C2").Value2 = "majidtv"

# ---------------------------------------------------------------------------
# 5) Sheet1 view: selection moves from B18 to B19 (dimension auto-extends).
# ---------------------------------------------------------------------------
$ws.Range("B19").Select() | Out-Null

# ---------------------------------------------------------------------------
# 6) Sheet2 (ImageLogic) view: drop the frozen/leftmost-column scroll
#    position (topLeftCell reverts back to A1 / default).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ImageLogic")
$ws2.Activate() | Out-Null
$ws2.Range("A1").Select() | Out-Null
$ws2.Range("I19").Select() | Out-Null

$ws.Activate() | Out-Null
